# Auto-generated edit script applying the Shinryu_Profits.xlsx commit diff
# Updates cached numeric values in ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 334.45456
$ws.Cells.Item(41, 10).Value = 666.6667
$ws.Cells.Item(41, 12).Value = 666.6667
$ws.Cells.Item(41, 14).Value = -1546.6667

$ws.Cells.Item(58, 8).Value = 1526.5
$ws.Cells.Item(58, 9).Value = 502.57144
$ws.Cells.Item(58, 10).Value = 2960
$ws.Cells.Item(58, 11).Value = 1507.71432
$ws.Cells.Item(58, 12).Value = 8880
$ws.Cells.Item(58, 13).Value = -1357.71432
$ws.Cells.Item(58, 14).Value = -9180

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1134.6129
$ws.Cells.Item(45, 9).Value = 850.2941
$ws.Cells.Item(45, 10).Value = 1479.8572
$ws.Cells.Item(45, 11).Value = 850.2941
$ws.Cells.Item(45, 12).Value = 1479.8572
$ws.Cells.Item(45, 13).Value = -473.2941
$ws.Cells.Item(45, 14).Value = -2233.8572

$ws.Cells.Item(97, 8).Value = 735.1905
$ws.Cells.Item(97, 9).Value = 390.6
$ws.Cells.Item(97, 10).Value = 1596.6666
$ws.Cells.Item(97, 11).Value = 390.6
$ws.Cells.Item(97, 12).Value = 1596.6666
$ws.Cells.Item(97, 13).Value = 105.4
$ws.Cells.Item(97, 14).Value = -2588.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 2674.5862
$ws.Cells.Item(94, 9).Value = 1227.6154
$ws.Cells.Item(94, 10).Value = 3850.25
$ws.Cells.Item(94, 11).Value = 1227.6154
$ws.Cells.Item(94, 12).Value = 3850.25
$ws.Cells.Item(94, 13).Value = -776.6153999999999
$ws.Cells.Item(94, 14).Value = -4752.25

$ws.Cells.Item(105, 8).Value = 913.1667
$ws.Cells.Item(105, 9).Value = 913.1667
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 913.1667
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 833.8333
$ws.Cells.Item(105, 14).ClearContents()

$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 3591.5833
$ws.Cells.Item(134, 9).Value = 1608.5
$ws.Cells.Item(134, 11).Value = 4825.5
$ws.Cells.Item(134, 13).Value = -2290.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 4130.7856
$ws.Cells.Item(3, 9).Value = 3214.5557
$ws.Cells.Item(3, 10).Value = 5780
$ws.Cells.Item(3, 11).Value = 9643.667099999999
$ws.Cells.Item(3, 12).Value = 17340
$ws.Cells.Item(3, 13).Value = -9531.667099999999
$ws.Cells.Item(3, 14).Value = -17564

$ws.Cells.Item(22, 8).Value = 1704.7273
$ws.Cells.Item(22, 10).Value = 1704.7273
$ws.Cells.Item(22, 12).Value = 5114.1819
$ws.Cells.Item(22, 14).Value = -5452.1819

$ws.Cells.Item(27, 8).Value = 1704.7273
$ws.Cells.Item(27, 10).Value = 1704.7273
$ws.Cells.Item(27, 12).Value = 5114.1819
$ws.Cells.Item(27, 14).Value = -5318.1819

$ws.Cells.Item(33, 8).Value = 96.833336
$ws.Cells.Item(33, 9).Value = 99.5
$ws.Cells.Item(33, 10).Value = 95.5
$ws.Cells.Item(33, 11).Value = 597
$ws.Cells.Item(33, 12).Value = 573
$ws.Cells.Item(33, 13).Value = -314
$ws.Cells.Item(33, 14).Value = -1139

$ws.Cells.Item(38, 8).Value = 73.86667
$ws.Cells.Item(38, 9).Value = 58
$ws.Cells.Item(38, 10).Value = 81.8
$ws.Cells.Item(38, 11).Value = 174
$ws.Cells.Item(38, 12).Value = 245.4
$ws.Cells.Item(38, 13).Value = 173
$ws.Cells.Item(38, 14).Value = -939.4

$ws.Cells.Item(68, 8).Value = 377.5
$ws.Cells.Item(68, 9).Value = 416
$ws.Cells.Item(68, 10).Value = 339
$ws.Cells.Item(68, 11).Value = 1248
$ws.Cells.Item(68, 12).Value = 1017
$ws.Cells.Item(68, 13).Value = -437
$ws.Cells.Item(68, 14).Value = -2639

$ws.Cells.Item(71, 8).Value = 377.5
$ws.Cells.Item(71, 9).Value = 416
$ws.Cells.Item(71, 10).Value = 339
$ws.Cells.Item(71, 11).Value = 3744
$ws.Cells.Item(71, 12).Value = 3051
$ws.Cells.Item(71, 13).Value = 312
$ws.Cells.Item(71, 14).Value = -11163

$ws.Cells.Item(86, 8).Value = 998.8261
$ws.Cells.Item(86, 9).Value = 788
$ws.Cells.Item(86, 10).Value = 1161
$ws.Cells.Item(86, 11).Value = 2364
$ws.Cells.Item(86, 12).Value = 3483
$ws.Cells.Item(86, 13).Value = -1178
$ws.Cells.Item(86, 14).Value = -5855

$ws.Cells.Item(89, 8).Value = 998.8261
$ws.Cells.Item(89, 9).Value = 788
$ws.Cells.Item(89, 10).Value = 1161
$ws.Cells.Item(89, 11).Value = 7092
$ws.Cells.Item(89, 12).Value = 10449
$ws.Cells.Item(89, 13).Value = -1164
$ws.Cells.Item(89, 14).Value = -22305

$ws.Cells.Item(136, 8).Value = 5036.8887
$ws.Cells.Item(136, 9).Value = 4466.3335
$ws.Cells.Item(136, 10).Value = 5322.1665
$ws.Cells.Item(136, 11).Value = 13399.0005
$ws.Cells.Item(136, 12).Value = 15966.4995
$ws.Cells.Item(136, 13).Value = -8299.000499999998
$ws.Cells.Item(136, 14).Value = -26166.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(43, 8).Value = 2746.2222
$ws.Cells.Item(43, 9).Value = 908.1429000000001
$ws.Cells.Item(43, 10).Value = 3915.9092
$ws.Cells.Item(43, 11).Value = 908.1429000000001
$ws.Cells.Item(43, 12).Value = 3915.9092
$ws.Cells.Item(43, 13).Value = -757.1429000000001
$ws.Cells.Item(43, 14).Value = -4217.9092

$ws.Cells.Item(97, 8).Value = 3620
$ws.Cells.Item(97, 9).Value = 4150
$ws.Cells.Item(97, 10).Value = 3266.6667
$ws.Cells.Item(97, 11).Value = 4150
$ws.Cells.Item(97, 12).Value = 3266.6667
$ws.Cells.Item(97, 13).Value = -3654
$ws.Cells.Item(97, 14).Value = -4258.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1395.6
$ws.Cells.Item(16, 9).Value = 500.625
$ws.Cells.Item(16, 10).Value = 4975.5
$ws.Cells.Item(16, 11).Value = 500.625
$ws.Cells.Item(16, 12).Value = 4975.5
$ws.Cells.Item(16, 13).Value = -330.625
$ws.Cells.Item(16, 14).Value = -5315.5

$ws.Cells.Item(46, 8).Value = 1037.5
$ws.Cells.Item(46, 9).Value = 900
$ws.Cells.Item(46, 10).Value = 2000
$ws.Cells.Item(46, 11).Value = 900
$ws.Cells.Item(46, 12).Value = 2000
$ws.Cells.Item(46, 13).Value = -712
$ws.Cells.Item(46, 14).Value = -2376

$ws.Cells.Item(55, 8).Value = 630.5714
$ws.Cells.Item(55, 9).Value = 595.2727
$ws.Cells.Item(55, 10).Value = 760
$ws.Cells.Item(55, 11).Value = 595.2727
$ws.Cells.Item(55, 12).Value = 760
$ws.Cells.Item(55, 13).Value = -422.2727
$ws.Cells.Item(55, 14).Value = -1106

$ws.Cells.Item(61, 8).Value = 1591.4
$ws.Cells.Item(61, 9).Value = 1238
$ws.Cells.Item(61, 11).Value = 1238
$ws.Cells.Item(61, 13).Value = -1036

$ws.Cells.Item(113, 8).Value = 1591.4
$ws.Cells.Item(113, 9).Value = 1238
$ws.Cells.Item(113, 11).Value = 1238
$ws.Cells.Item(113, 13).Value = 932

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 772.4286
$ws.Cells.Item(100, 9).Value = 488
$ws.Cells.Item(100, 10).Value = 1151.6666
$ws.Cells.Item(100, 11).Value = 976
$ws.Cells.Item(100, 12).Value = 2303.3332
$ws.Cells.Item(100, 13).Value = -435
$ws.Cells.Item(100, 14).Value = -3385.3332
